$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.141.54"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").Value = "3.454.79"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.90"
$ws.Range("E5").Value = "  +3.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.88"
$ws.Range("E6").Value = "  +6.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "3.446.62"
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.59"
$ws.Range("E12").Value = "  +6.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000277"
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.48"
$ws.Range("E14").Value = "  +2.86%  "
$ws.Range("D15").Value = "3.999.01"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.94"
$ws.Range("E16").Value = "  +2.82%  "
$ws.Range("D17").Value = "3.454.67"
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("D18").Value = "67.237.54"
$ws.Range("E18").Value = "  +2.49%  "
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.05"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "481.69"
$ws.Range("E22").Value = "  +5.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.36"
$ws.Range("E23").Value = "  +23.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.35"
$ws.Range("E24").Value = "  +9.00%  "
$ws.Range("E25").Value = "  +5.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.41"
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.98"
$ws.Range("E28").Value = "  +2.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.01"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.25"
$ws.Range("E30").Value = "  +0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.38"
$ws.Range("E31").Value = "  +12.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "604.68"
$ws.Range("E32").Value = "  +3.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.86"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.81"
$ws.Range("E34").Value = "  +2.52%  "
$ws.Range("E35").Value = "  +3.63%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.146"
$ws.Range("E37").Value = "  +2.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.99"
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.387"
$ws.Range("E39").Value = "  +3.53%  "
$ws.Range("D40").Value = "0.0₃0770"
$ws.Range("E40").Value = "  +3.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.47"
$ws.Range("E41").Value = "  -3.94%  "
$ws.Range("D42").Value = "3.201.47"
$ws.Range("E42").Value = "  +3.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.88"
$ws.Range("E43").Value = "  +3.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0430"
$ws.Range("E44").Value = "  +2.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.57"
$ws.Range("E45").Value = "  +4.75%  "
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.69"
$ws.Range("E48").Value = "  +16.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.66"
$ws.Range("E50").Value = "  +3.34%  "
$ws.Range("E51").Value = "  +2.63%  "
